$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) matching the header formatting of the
# existing header row (reuse G1's cell format so no new style is minted).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the Save flag for each data row (1 = save recorded that outing).
$saveValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 1
    8 = 1
    9 = 0
    10 = 0
    11 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
